$d = $word.ActiveDocument

# 1) Notary line: hardcoded "Texas" becomes the client's actual address state field.
$d.Content.Find.Execute(
    "Notary Public, State of Texas",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Notary Public, State of {{doc.case.client[0].address.state}}",
    2
)

# 2) "Executed in" line: county/state now come from the client's address instead of the case.
$d.Content.Find.Execute(
    "Executed in {{doc.case.county}} County, {{doc.case.state}}, on",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Executed in {{doc.case.client[0].address.county}}, {{doc.case.client[0].address.state}}, on",
    2
)
